$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2346.3057
$ws.Range("I15").Value = 2346.3057
$ws.Range("K15").Value = 7038.9171
$ws.Range("M15").Value = -6869.9171
$ws.Range("H41").Value = 311.58334
$ws.Range("I41").Value = 293.9
$ws.Range("K41").Value = 293.9
$ws.Range("M41").Value = 146.1
$ws.Range("I62").Value = 11259.667
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 11259.667
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -10635.667
$ws.Range("N62").Value = -16248
$ws.Range("H64").Value = 10154.866
$ws.Range("I64").Value = 7173
$ws.Range("J64").Value = 10900.333
$ws.Range("K64").Value = 7173
$ws.Range("L64").Value = 10900.333
$ws.Range("M64").Value = -6925
$ws.Range("N64").Value = -11396.333
$ws.Range("I65").Value = 11259.667
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 56298.335
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -53178.335
$ws.Range("N65").Value = -81240
$ws.Range("H67").Value = 10154.866
$ws.Range("I67").Value = 7173
$ws.Range("J67").Value = 10900.333
$ws.Range("K67").Value = 7173
$ws.Range("L67").Value = 10900.333
$ws.Range("M67").Value = -6315
$ws.Range("N67").Value = -12616.333
$ws.Range("H103").Value = 454
$ws.Range("I103").Value = 469.5
$ws.Range("J103").Value = 433.33334
$ws.Range("K103").Value = 1408.5
$ws.Range("L103").Value = 1300.00002
$ws.Range("M103").Value = -822.5
$ws.Range("N103").Value = -2472.00002
$ws.Range("H111").Value = 6152
$ws.Range("I111").Value = 4228.5
$ws.Range("K111").Value = 12685.5
$ws.Range("M111").Value = -9618.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4804.1875
$ws.Range("I32").Value = 4155.738
$ws.Range("K32").Value = 4155.738
$ws.Range("M32").Value = -3868.738
$ws.Range("H45").Value = 159567.08
$ws.Range("I45").Value = 337066.5
$ws.Range("K45").Value = 337066.5
$ws.Range("M45").Value = -336689.5
$ws.Range("H102").Value = 3016.9583
$ws.Range("I102").Value = 2290.35
$ws.Range("K102").Value = 2290.35
$ws.Range("M102").Value = -668.3499999999999
$ws.Range("H132").Value = 4787.5884
$ws.Range("I132").Value = 4143.5557
$ws.Range("K132").Value = 12430.6671
$ws.Range("M132").Value = -9900.667099999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3854.1
$ws.Range("I20").Value = 5153.4
$ws.Range("J20").Value = 2554.8
$ws.Range("K20").Value = 5153.4
$ws.Range("L20").Value = 2554.8
$ws.Range("M20").Value = -4906.4
$ws.Range("N20").Value = -3048.8

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3234.4902
$ws.Range("I31").Value = 2070.2163
$ws.Range("K31").Value = 2070.2163
$ws.Range("M31").Value = -1775.2163
$ws.Range("H34").Value = 3234.4902
$ws.Range("I34").Value = 2070.2163
$ws.Range("K34").Value = 2070.2163
$ws.Range("M34").Value = -1868.2163
$ws.Range("H58").Value = 5539.8125
$ws.Range("I58").Value = 2649.6667
$ws.Range("J58").Value = 9255.714
$ws.Range("K58").Value = 2649.6667
$ws.Range("L58").Value = 9255.714
$ws.Range("M58").Value = -2446.6667
$ws.Range("N58").Value = -9661.714
$ws.Range("H136").Value = 5539.8125
$ws.Range("I136").Value = 2649.6667
$ws.Range("J136").Value = 9255.714
$ws.Range("K136").Value = 7949.000100000001
$ws.Range("L136").Value = 27767.142
$ws.Range("M136").Value = -5399.000100000001
$ws.Range("N136").Value = -32867.142

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1425.9429
$ws.Range("I5").Value = 1429.9445
$ws.Range("J5").Value = 1421.7059
$ws.Range("K5").Value = 4289.833500000001
$ws.Range("L5").Value = 4265.1177
$ws.Range("M5").Value = -4177.833500000001
$ws.Range("N5").Value = -4489.1177
$ws.Range("H104").Value = 4266.6665
$ws.Range("I104").Value = 4800
$ws.Range("J104").Value = 4000
$ws.Range("K104").Value = 14400
$ws.Range("L104").Value = 12000
$ws.Range("M104").Value = -11779
$ws.Range("N104").Value = -17242
$ws.Range("H114").Value = 2917.1428
$ws.Range("I114").Value = 3000
$ws.Range("J114").Value = 2903.3333
$ws.Range("K114").Value = 9000
$ws.Range("L114").Value = 8709.999899999999
$ws.Range("M114").Value = -5746
$ws.Range("N114").Value = -15217.9999
$ws.Range("H131").Value = 38463868
$ws.Range("I131").Value = 166667140
$ws.Range("J131").Value = 2883.5
$ws.Range("K131").Value = 500001420
$ws.Range("L131").Value = 8650.5
$ws.Range("M131").Value = -499996380
$ws.Range("N131").Value = -18730.5
$ws.Range("H135").Value = 1425.9429
$ws.Range("I135").Value = 1429.9445
$ws.Range("J135").Value = 1421.7059
$ws.Range("K135").Value = 12869.5005
$ws.Range("L135").Value = 12795.3531
$ws.Range("M135").Value = -10334.5005
$ws.Range("N135").Value = -17865.3531

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14226.125
$ws.Range("J70").Value = 15815.546
$ws.Range("L70").Value = 15815.546
$ws.Range("N70").Value = -16355.546
$ws.Range("H73").Value = 14226.125
$ws.Range("J73").Value = 15815.546
$ws.Range("L73").Value = 15815.546
$ws.Range("N73").Value = -17687.546
$ws.Range("H102").Value = 3510.4285
$ws.Range("I102").Value = 2914.7
$ws.Range("K102").Value = 2914.7
$ws.Range("M102").Value = -1292.7

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7058.52
$ws.Range("I7").Value = 8400.857
$ws.Range("J7").Value = 5350.091
$ws.Range("K7").Value = 8400.857
$ws.Range("L7").Value = 5350.091
$ws.Range("M7").Value = -8288.857
$ws.Range("N7").Value = -5574.091
$ws.Range("H40").Value = 4269
$ws.Range("I40").Value = 3856.7144
$ws.Range("K40").Value = 3856.7144
$ws.Range("M40").Value = -3720.7144
$ws.Range("H107").Value = 3033.3333
$ws.Range("I107").Value = 3033.3333
$ws.Range("K107").Value = 3033.3333
$ws.Range("M107").Value = -1113.3333
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
$ws.Range("H126").Value = 7058.52
$ws.Range("I126").Value = 8400.857
$ws.Range("J126").Value = 5350.091
$ws.Range("K126").Value = 25202.571
$ws.Range("L126").Value = 16050.273
$ws.Range("M126").Value = -22732.571
$ws.Range("N126").Value = -20990.273
$ws.Range("H132").Value = 25001.6
$ws.Range("I132").Value = 56504
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 169512
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -166982
$ws.Range("N132").Value = -17060

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 811.1111
$ws.Range("I107").Value = 811.1111
$ws.Range("K107").Value = 2433.3333
$ws.Range("M107").Value = -513.3332999999998
$ws.Range("H109").Value = 69000
$ws.Range("I109").Value = 69000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 69000
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -67613
$ws.Range("N109").ClearContents()
$ws.Range("H116").Value = 132500
$ws.Range("J116").Value = 132500
$ws.Range("L116").Value = 132500
$ws.Range("N116").Value = -141678
$ws.Range("H122").Value = 4095.9285
$ws.Range("I122").Value = 3461.375
$ws.Range("K122").Value = 10384.125
$ws.Range("M122").Value = -7934.125
$ws.Range("H132").Value = 4218.6587
$ws.Range("I132").Value = 3569.5518
$ws.Range("K132").Value = 10708.6554
$ws.Range("M132").Value = -8178.6554
$ws.Range("H136").Value = 5053.8965
$ws.Range("I136").Value = 4056.6316
$ws.Range("J136").Value = 6948.7
$ws.Range("K136").Value = 12169.8948
$ws.Range("L136").Value = 20846.1
$ws.Range("M136").Value = -9619.8948
$ws.Range("N136").Value = -25946.1
